$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H34").Value = 10500
$ws.Range("I34").Value = 1000
$ws.Range("K34").Value = 1000
$ws.Range("M34").Value = -797

$ws.Range("H36").Value = 10500
$ws.Range("I36").Value = 1000
$ws.Range("K36").Value = 1000
$ws.Range("M36").Value = -285

$ws.Range("H74").Value = 3240.2666
$ws.Range("I74").Value = 2885.4
$ws.Range("K74").Value = 2885.4
$ws.Range("M74").Value = -1949.4

$ws.Range("H76").Value = 2896.26
$ws.Range("I76").Value = 2815.4243
$ws.Range("J76").Value = 3053.1765
$ws.Range("K76").Value = 2815.4243
$ws.Range("L76").Value = 3053.1765
$ws.Range("M76").Value = -2500.4243
$ws.Range("N76").Value = -3683.1765

$ws.Range("H77").Value = 3240.2666
$ws.Range("I77").Value = 2885.4
$ws.Range("K77").Value = 14427
$ws.Range("M77").Value = -9747

$ws.Range("H79").Value = 2896.26
$ws.Range("I79").Value = 2815.4243
$ws.Range("J79").Value = 3053.1765
$ws.Range("K79").Value = 2815.4243
$ws.Range("L79").Value = 3053.1765
$ws.Range("M79").Value = -1723.4243
$ws.Range("N79").Value = -5237.1765

$ws.Range("H127").Value = 2154736.5
$ws.Range("I127").Value = 497
$ws.Range("K127").Value = 1491
$ws.Range("M127").Value = 3469

$ws.Range("H138").Value = 15878090
$ws.Range("I138").Value = 22226260
$ws.Range("J138").Value = 7666.6665
$ws.Range("K138").Value = 66678780
$ws.Range("L138").Value = 22999.9995
$ws.Range("M138").Value = -66673640
$ws.Range("N138").Value = -33279.99950000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 21357.475
$ws.Range("I32").Value = 21292.428
$ws.Range("J32").Value = 25000
$ws.Range("K32").Value = 21292.428
$ws.Range("L32").Value = 25000
$ws.Range("M32").Value = -21005.428
$ws.Range("N32").Value = -25574

$ws.Range("H55").Value = 33168.668
$ws.Range("J55").Value = 39202.4
$ws.Range("L55").Value = 39202.4
$ws.Range("N55").Value = -39832.4

$ws.Range("H61").Value = 1813.3485
$ws.Range("I61").Value = 1885.9
$ws.Range("J61").Value = 1586.625
$ws.Range("K61").Value = 1885.9
$ws.Range("L61").Value = 1586.625
$ws.Range("M61").Value = -1673.9
$ws.Range("N61").Value = -2010.625

$ws.Range("H97").Value = 5987.5264
$ws.Range("I97").Value = 6726.875
$ws.Range("J97").Value = 2044.3334
$ws.Range("K97").Value = 6726.875
$ws.Range("L97").Value = 2044.3334
$ws.Range("M97").Value = -6230.875
$ws.Range("N97").Value = -3036.3334

$ws.Range("H102").Value = 1514.2667
$ws.Range("I102").Value = 1234.8334
$ws.Range("J102").Value = 2632
$ws.Range("K102").Value = 1234.8334
$ws.Range("L102").Value = 2632
$ws.Range("M102").Value = 387.1666
$ws.Range("N102").Value = -5876

$ws.Range("H110").Value = 524.05884
$ws.Range("I110").Value = 499.6
$ws.Range("K110").Value = 499.6
$ws.Range("M110").Value = 1545.4

$ws.Range("H136").Value = 1813.3485
$ws.Range("I136").Value = 1885.9
$ws.Range("J136").Value = 1586.625
$ws.Range("K136").Value = 5657.700000000001
$ws.Range("L136").Value = 4759.875
$ws.Range("M136").Value = -3107.700000000001
$ws.Range("N136").Value = -9859.875

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 3264.3264
$ws.Range("I105").Value = 1942.7407
$ws.Range("J105").Value = 4886.273
$ws.Range("K105").Value = 1942.7407
$ws.Range("L105").Value = 4886.273
$ws.Range("M105").Value = -195.7407000000001
$ws.Range("N105").Value = -8380.273000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6070.75
$ws.Range("I31").Value = 2019.7407
$ws.Range("J31").Value = 18223.777
$ws.Range("K31").Value = 2019.7407
$ws.Range("L31").Value = 18223.777
$ws.Range("M31").Value = -1724.7407
$ws.Range("N31").Value = -18813.777

$ws.Range("H34").Value = 6070.75
$ws.Range("I34").Value = 2019.7407
$ws.Range("J34").Value = 18223.777
$ws.Range("K34").Value = 2019.7407
$ws.Range("L34").Value = 18223.777
$ws.Range("M34").Value = -1817.7407
$ws.Range("N34").Value = -18627.777

$ws.Range("H36").Value = 5000
$ws.Range("I36").Value = 5000
$ws.Range("K36").Value = 5000
$ws.Range("M36").Value = -4612

$ws.Range("H40").Value = 5000
$ws.Range("I40").Value = 5000
$ws.Range("K40").Value = 5000
$ws.Range("M40").Value = -4840

$ws.Range("H75").Value = 43185
$ws.Range("J75").Value = 43185
$ws.Range("L75").Value = 43185
$ws.Range("N75").Value = -45181

$ws.Range("H78").Value = 43185
$ws.Range("J78").Value = 43185
$ws.Range("L78").Value = 129555
$ws.Range("N78").Value = -139539

$ws.Range("H134").Value = 1702703.9
$ws.Range("I134").Value = 2679.75
$ws.Range("J134").Value = 3969402.8
$ws.Range("K134").Value = 8039.25
$ws.Range("L134").Value = 11908208.4
$ws.Range("M134").Value = -5504.25
$ws.Range("N134").Value = -11913278.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 661441.7
$ws.Range("I2").Value = 52.875
$ws.Range("J2").Value = 1068450.1
$ws.Range("K2").Value = 317.25
$ws.Range("L2").Value = 6410700.600000001
$ws.Range("M2").Value = -204.25
$ws.Range("N2").Value = -6410926.600000001

$ws.Range("H17").Value = 50
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 50
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 150
$ws.Range("M17").Value = $null
$ws.Range("N17").Value = -488

$ws.Range("H34").Value = 241.42857
$ws.Range("I34").Value = 115
$ws.Range("J34").Value = 1000
$ws.Range("K34").Value = 345
$ws.Range("L34").Value = 3000
$ws.Range("M34").Value = -261
$ws.Range("N34").Value = -3168

$ws.Range("H39").Value = 496.55173
$ws.Range("I39").Value = 450
$ws.Range("J39").Value = 500
$ws.Range("K39").Value = 1350
$ws.Range("L39").Value = 1500
$ws.Range("M39").Value = -1056
$ws.Range("N39").Value = -2088

$ws.Range("H55").Value = 439.66666
$ws.Range("J55").Value = 550
$ws.Range("L55").Value = 1650
$ws.Range("N55").Value = -2004

$ws.Range("H107").Value = 423.9189
$ws.Range("I107").Value = 196.53847
$ws.Range("J107").Value = 547.0833
$ws.Range("K107").Value = 589.61541
$ws.Range("L107").Value = 1641.2499
$ws.Range("M107").Value = 1330.38459
$ws.Range("N107").Value = -5481.2499

$ws.Range("H113").Value = 899.04083
$ws.Range("I113").Value = 681.9091
$ws.Range("J113").Value = 1346.875
$ws.Range("K113").Value = 2045.7273
$ws.Range("L113").Value = 4040.625
$ws.Range("M113").Value = 124.2727
$ws.Range("N113").Value = -8380.625

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 9000
$ws.Range("I5").Value = 9000
$ws.Range("K5").Value = 9000
$ws.Range("M5").Value = -8888

$ws.Range("H70").Value = 9179.684999999999
$ws.Range("I70").Value = 12230.5
$ws.Range("J70").Value = 3949.7144
$ws.Range("K70").Value = 12230.5
$ws.Range("L70").Value = 3949.7144
$ws.Range("M70").Value = -11960.5
$ws.Range("N70").Value = -4489.7144

$ws.Range("H73").Value = 9179.684999999999
$ws.Range("I73").Value = 12230.5
$ws.Range("J73").Value = 3949.7144
$ws.Range("K73").Value = 12230.5
$ws.Range("L73").Value = 3949.7144
$ws.Range("M73").Value = -11294.5
$ws.Range("N73").Value = -5821.7144

$ws.Range("H113").Value = 67881.8
$ws.Range("I113").Value = 84516.75
$ws.Range("J113").Value = 1342
$ws.Range("K113").Value = 84516.75
$ws.Range("L113").Value = 1342
$ws.Range("M113").Value = -82346.75
$ws.Range("N113").Value = -5682

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H48").Value = 14560
$ws.Range("I48").Value = 13300
$ws.Range("J48").Value = 16450
$ws.Range("K48").Value = 13300
$ws.Range("L48").Value = 16450
$ws.Range("M48").Value = -12639
$ws.Range("N48").Value = -17772

$ws.Range("H56").Value = 15458.5
$ws.Range("J56").Value = 21514.166
$ws.Range("L56").Value = 21514.166
$ws.Range("N56").Value = -22896.166

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H51").Value = 13318
$ws.Range("I51").Value = 3800
$ws.Range("K51").Value = 3800
$ws.Range("M51").Value = -3290

$ws.Range("H132").Value = 1774.5143
$ws.Range("I132").Value = 1631.579
$ws.Range("J132").Value = 1944.25
$ws.Range("K132").Value = 4894.737
$ws.Range("L132").Value = 5832.75
$ws.Range("M132").Value = -2364.737
$ws.Range("N132").Value = -10892.75
